$d = $word.ActiveDocument

$p4 = $d.Paragraphs(4)
$target = $d.Range($p4.Range.Start, $p4.Range.End)

$inner = '<w:p><w:pPr><w:rPr><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>All Tournaments Page:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Fixed layout alignment in all tournaments page.</w:t></w:r><w:r><w:rPr><w:bCs/><w:color w:val="EE0000"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:bCs/></w:rPr><w:t>fixed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:bCs/></w:rPr><w:t xml:space="preserve"> up the visit alignment in favourite tournaments</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:br/><w:t>added tooltips to remaining icons.</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:br/><w:t>fixed up choose file button with new label button to for styling</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:bCs/><w:color w:val="C00000"/></w:rPr><w:t>needs reset or delete button</w:t></w:r><w:r><w:rPr><w:bCs/><w:color w:val="C00000"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:bCs/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>HeaderBar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:br/></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>star as settings is confusing</w:t></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:bCs/></w:rPr><w:t>gear, sandwich bar</w:t></w:r></w:p>'

$xml = "<?xml version='1.0' standalone='yes'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>" + $inner + "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

$target.InsertXML($xml)
